$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 7656.3335
$ws.Range("J7").Value = 7656.3335
$ws.Range("L7").Value = 7656.3335
$ws.Range("N7").Value = -7880.3335
$ws.Range("H14").Value = 7656.3335
$ws.Range("J14").Value = 7656.3335
$ws.Range("L14").Value = 7656.3335
$ws.Range("N14").Value = -8038.3335
$ws.Range("H40").Value = 4938.5
$ws.Range("I40").Value = 4250.25
$ws.Range("J40").Value = 6315
$ws.Range("K40").Value = 4250.25
$ws.Range("L40").Value = 6315
$ws.Range("M40").Value = -4075.25
$ws.Range("N40").Value = -6665
$ws.Range("H92").Value = 741.0909
$ws.Range("I92").Value = 735.2
$ws.Range("K92").Value = 735.2
$ws.Range("M92").Value = 512.8
$ws.Range("H106").Value = 3441.6667
$ws.Range("I106").Value = 3325
$ws.Range("K106").Value = 3325
$ws.Range("M106").Value = -2694
$ws.Range("H135").Value = 2541.6667
$ws.Range("I135").Value = 2443.5
$ws.Range("K135").Value = 21991.5
$ws.Range("M135").Value = -19456.5
$ws.Range("H137").Value = 10218.454
$ws.Range("I137").Value = 4165.3335
$ws.Range("J137").Value = 12488.375
$ws.Range("K137").Value = 12496.0005
$ws.Range("L137").Value = 37465.125
$ws.Range("M137").Value = -9946.000499999998
$ws.Range("N137").Value = -42565.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 2925
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5346
$ws.Range("H61").Value = 7854.1665
$ws.Range("I61").Value = 6739.7437
$ws.Range("K61").Value = 6739.7437
$ws.Range("M61").Value = -6527.7437
$ws.Range("H74").Value = 457635.47
$ws.Range("I74").Value = 626893.9
$ws.Range("K74").Value = 626893.9
$ws.Range("M74").Value = -626019.9
$ws.Range("H77").Value = 457635.47
$ws.Range("I77").Value = 626893.9
$ws.Range("K77").Value = 3134469.5
$ws.Range("M77").Value = -3130101.5
$ws.Range("H132").Value = 3377.2834
$ws.Range("I132").Value = 2948.6492
$ws.Range("K132").Value = 8845.9476
$ws.Range("M132").Value = -6315.9476
$ws.Range("H136").Value = 7854.1665
$ws.Range("I136").Value = 6739.7437
$ws.Range("K136").Value = 20219.2311
$ws.Range("M136").Value = -17669.2311

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 3867
$ws.Range("I22").Value = 3867
$ws.Range("K22").Value = 3867
$ws.Range("M22").Value = -3694
$ws.Range("H29").Value = 1488.3334
$ws.Range("I29").Value = 1488.3334
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 1488.3334
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -1199.3334
$ws.Range("H94").Value = 1796.5555
$ws.Range("I94").Value = 1566.2727
$ws.Range("K94").Value = 1566.2727
$ws.Range("M94").Value = -1115.2727

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 5124.75
$ws.Range("I33").Value = 1250
$ws.Range("J33").Value = 8999.5
$ws.Range("K33").Value = 1250
$ws.Range("L33").Value = 8999.5
$ws.Range("M33").Value = -871
$ws.Range("N33").Value = -9757.5
$ws.Range("H132").Value = 3369.25
$ws.Range("I132").Value = 1756
$ws.Range("K132").Value = 5268
$ws.Range("M132").Value = -2738

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1920536.5
$ws.Range("I4").Value = 3054765.2
$ws.Range("J4").Value = 1072.3846
$ws.Range("K4").Value = 9164295.600000001
$ws.Range("L4").Value = 3217.1538
$ws.Range("M4").Value = -9164183.600000001
$ws.Range("N4").Value = -3441.1538
$ws.Range("H113").Value = 653.5454999999999
$ws.Range("J113").Value = 1084.75
$ws.Range("L113").Value = 3254.25
$ws.Range("N113").Value = -7594.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1681122.5
$ws.Range("J7").Value = 1444443.1
$ws.Range("L7").Value = 1444443.1
$ws.Range("N7").Value = -1444667.1
$ws.Range("H8").Value = 1681122.5
$ws.Range("J8").Value = 1444443.1
$ws.Range("L8").Value = 1444443.1
$ws.Range("N8").Value = -1444721.1
$ws.Range("H9").Value = 11156.333
$ws.Range("I9").Value = 13250
$ws.Range("J9").Value = 6969
$ws.Range("K9").Value = 13250
$ws.Range("L9").Value = 6969
$ws.Range("M9").Value = -13080
$ws.Range("N9").Value = -7309
$ws.Range("H113").Value = 16749.412
$ws.Range("I113").Value = 17758.75
$ws.Range("J113").Value = 600
$ws.Range("K113").Value = 17758.75
$ws.Range("L113").Value = 600
$ws.Range("M113").Value = -15588.75
$ws.Range("N113").Value = -4940
$ws.Range("H132").Value = 5644.2334
$ws.Range("I132").Value = 3934.5833
$ws.Range("K132").Value = 11803.7499
$ws.Range("M132").Value = -9273.749899999999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 60040.332
$ws.Range("I45").Value = 60040.332
$ws.Range("K45").Value = 60040.332
$ws.Range("M45").Value = -59633.332
$ws.Range("H46").Value = 7770.0835
$ws.Range("I46").Value = 1167
$ws.Range("J46").Value = 8713.380999999999
$ws.Range("K46").Value = 1167
$ws.Range("L46").Value = 8713.380999999999
$ws.Range("M46").Value = -979
$ws.Range("N46").Value = -9089.380999999999
$ws.Range("H61").Value = 1167.8636
$ws.Range("I61").Value = 1004.6316
$ws.Range("J61").Value = 2201.6667
$ws.Range("K61").Value = 1004.6316
$ws.Range("L61").Value = 2201.6667
$ws.Range("M61").Value = -802.6316
$ws.Range("N61").Value = -2605.6667
$ws.Range("H113").Value = 1167.8636
$ws.Range("I113").Value = 1004.6316
$ws.Range("J113").Value = 2201.6667
$ws.Range("K113").Value = 1004.6316
$ws.Range("L113").Value = 2201.6667
$ws.Range("M113").Value = 1165.3684
$ws.Range("N113").Value = -6541.6667
$ws.Range("H141").Value = 99999.75
$ws.Range("J141").Value = 99999.75
$ws.Range("L141").Value = 99999.75
$ws.Range("N141").Value = -110359.75

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30000
$ws.Range("I40").Value = 30000
$ws.Range("K40").Value = 30000
$ws.Range("M40").Value = -29851
$ws.Range("H81").Value = 11852.283
$ws.Range("I81").Value = 2798.7334
$ws.Range("J81").Value = 15426.053
$ws.Range("K81").Value = 5597.4668
$ws.Range("L81").Value = 30852.106
$ws.Range("M81").Value = -4536.4668
$ws.Range("N81").Value = -32974.106
$ws.Range("H84").Value = 11852.283
$ws.Range("I84").Value = 2798.7334
$ws.Range("J84").Value = 15426.053
$ws.Range("K84").Value = 27987.334
$ws.Range("L84").Value = 154260.53
$ws.Range("M84").Value = -22683.334
$ws.Range("N84").Value = -164868.53
$ws.Range("H122").Value = 2481.7368
$ws.Range("I122").Value = 2508.5
$ws.Range("K122").Value = 7525.5
$ws.Range("M122").Value = -5075.5
$ws.Range("H132").Value = 5473.8184
$ws.Range("I132").Value = 4716.5435
$ws.Range("J132").Value = 9344.333000000001
$ws.Range("K132").Value = 14149.6305
$ws.Range("L132").Value = 28032.999
$ws.Range("M132").Value = -11619.6305
$ws.Range("N132").Value = -33092.999
$ws.Range("H136").Value = 2782.4
$ws.Range("I136").Value = 1249.6086
$ws.Range("K136").Value = 3748.8258
$ws.Range("M136").Value = -1198.8258
